{"js": "// Q&A placeholder cleanup:\n//   Para 1 (runs \"S\" + \"dfasdfsdfsdf\" + \"dsfgsdg...\")      -> \"What is 3+3?\"\n//   Para 2 (\"F\") + Para 3 (\"Asdf\") + Para 4 (empty)\n//     + Para 5 (\"Asdflkjsda;kfljasdfkl\", holds _GoBack bookmark)\n//   collapse into a single paragraph: \"Ans:\" + \" 7\", bookmark kept.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraphs 3 & 4 (0-based indices 2 and 3 -- \"Asdf\" and the blank\n// paragraph) disappear entirely; queue both deletes before syncing.\nparagraphs.items[2].delete();\nparagraphs.items[3].delete();\nawait context.sync();\n\n// Re-fetch: the body now has 3 paragraphs -- \"S.../F/Asdflkjsda...\".\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\n\nconst firstPara = remaining.items[0];   // \"S\" + \"dfasdfsdfsdf\" + \"dsfgsdg...\"\nconst fParaRange = remaining.items[1].getRange(\"End\");        // end of \"F\" paragraph\nconst bookmarkParaStart = remaining.items[2].getRange(\"Start\"); // start of \"Asdflkjsda...\" paragraph\n\n// Deleting the span between them removes the intervening paragraph mark,\n// merging \"F\" and \"Asdflkjsda;kfljasdfkl\" (plus its bookmark) into one\n// paragraph: \"FAsdflkjsda;kfljasdfkl\".\nconst paraMark = fParaRange.expandTo(bookmarkParaStart);\nparaMark.delete();\nawait context.sync();\n\n// Replace paragraph 1's whole content with the question.\nfirstPara.getRange(\"Whole\").insertText(\"What is 3+3?\", \"Replace\");\nawait context.sync();\n\n// The merged second paragraph still has two distinct runs: \"F\" and\n// \"Asdflkjsda;kfljasdfkl\". Retarget each in place (so the run split /\n// bookmark position are preserved) instead of rewriting the paragraph.\nconst mergedParagraphs = body.paragraphs;\nmergedParagraphs.load(\"items\");\nawait context.sync();\n\nconst answerPara = mergedParagraphs.items[1];\nconst answerRange = answerPara.getRange(\"Whole\");\n\nconst fHits = answerRange.search(\"F\", { matchCase: true });\nfHits.load(\"items\");\nawait context.sync();\nfHits.items[0].insertText(\"Ans:\", \"Replace\");\nawait context.sync();\n\nconst tailHits = answerPara\n  .getRange(\"Whole\")\n  .search(\"Asdflkjsda;kfljasdfkl\", { matchCase: true });\ntailHits.load(\"items\");\nawait context.sync();\ntailHits.items[0].insertText(\" 7\", \"Replace\");\nawait context.sync();\n", "ps1": "# Q&A placeholder cleanup:\n#   Para 1 (runs \"S\" + \"dfasdfsdfsdf\" + \"dsfgsdg...\")      -> \"What is 3+3?\"\n#   Para 2 (\"F\") + Para 3 (\"Asdf\") + Para 4 (empty)\n#     + Para 5 (\"Asdflkjsda;kfljasdfkl\", holds the _GoBack bookmark)\n#   collapse into a single paragraph: \"Ans:\" + \" 7\", bookmark kept.\n\n$d = $word.ActiveDocument\n\n# Paragraphs 3 & 4 (\"Asdf\" and the blank paragraph) disappear entirely.\n# After the first delete, what was paragraph 4 shifts down to index 3, so\n# deleting index 3 twice removes both.\n$d.Paragraphs(3).Range.Delete()\n$d.Paragraphs(3).Range.Delete()\n\n# Document is now 3 paragraphs:\n#   1: \"S\" + \"dfasdfsdfsdf\" + \"dsfgsdg...\"\n#   2: \"F\"\n#   3: \"Asdflkjsda;kfljasdfkl\" (+ the _GoBack bookmark)\n\n# Rewrite paragraph 1's whole content (collapses its 3 runs into one).\n$p1 = $d.Paragraphs(1)\n$r1 = $p1.Range\n$r1.End = $r1.End - 1          # exclude the paragraph mark from the replace\n$r1.Text = \"What is 3+3?\"\n\n# Rewrite paragraph 2's text in place, before merging -- this keeps the\n# edit scoped to a single run instead of touching the bookmark paragraph.\n$p2 = $d.Paragraphs(2)\n$r2 = $p2.Range\n$r2.End = $r2.End - 1\n$r2.Text = \"Ans:\"\n\n# Rewrite paragraph 3's text in place too (bookmark untouched since the\n# replace range stops before the paragraph mark that follows it).\n$p3 = $d.Paragraphs(3)\n$r3 = $p3.Range\n$r3.End = $r3.End - 1\n$r3.Text = \" 7\"\n\n# Finally, merge paragraph 2 (\"Ans:\") with paragraph 3 (\" 7\" + bookmark)\n# by deleting the paragraph mark between them -- same as pressing Delete\n# at the end of line 2. The two runs stay distinct and the bookmark rides\n# along into the merged paragraph.\n$p2 = $d.Paragraphs(2)\n$p3 = $d.Paragraphs(3)\n$joinRange = $d.Range($p2.Range.End - 1, $p3.Range.Start)\n$joinRange.Delete()\n"}
